$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F6").Value = 1321
$ws1.Range("F10").Value = 448
$ws1.Range("F12").Value = 214
$ws1.Range("F15").Value = 471
$ws1.Range("F17").Value = 1053
$ws1.Range("F19").Value = 290
$ws1.Range("F20").Value = 413
$ws1.Range("F26").Value = 447
$ws1.Range("F27").Value = 298

$ws2.Range("F4").Value = 381
$ws2.Range("F5").Value = 50
$ws2.Range("F6").Value = 49

$ws4.Range("F8").Value = 1321
$ws4.Range("F11").Value = 381
$ws4.Range("F13").Value = 50
$ws4.Range("F15").Value = 49
$ws4.Range("F17").Value = 448
$ws4.Range("F19").Value = 214
$ws4.Range("F22").Value = 471
$ws4.Range("F24").Value = 1053
$ws4.Range("F28").Value = 290
$ws4.Range("F29").Value = 413
$ws4.Range("F41").Value = 447
$ws4.Range("F42").Value = 298
